$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.403.99'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '2.018.41'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '260.48'
$ws.Range("E5").Value = '  +5.87%  '
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '56.54'
$ws.Range("E8").Value = '  -5.65%  '
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = '0.0775'
$ws.Range("E10").Value = '  -3.79%  '
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("E12").Value = '  -4.39%  '
$ws.Range("D13").Value = '2.315.10'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("E14").Value = '  -4.01%  '
$ws.Range("D15").Value = '21.03'
$ws.Range("E15").Value = '  -6.26%  '
$ws.Range("E16").Value = '  -2.71%  '
$ws.Range("D17").Value = '2.032.32'
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("D18").Value = '37.386.33'
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '70.01'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("E20").Value = '  -2.64%  '
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = '227.72'
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").Value = '2.65'
$ws.Range("E23").Value = '  +7.65%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").Value = '165.08'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("E27").Value = '  -4.49%  '
$ws.Range("D28").Value = '19.77'
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("E29").Value = '  -6.29%  '
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").Value = '4.59'
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("E35").Value = '  +0.63%  '
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").Value = '3.37'
$ws.Range("E37").Value = '  +1.44%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").Value = '5.28'
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("E40").Value = '  +4.02%  '
$ws.Range("E41").Value = '  +2.23%  '
$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").Value = '0.0938'
$ws.Range("E42").Value = '  -4.81%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0214'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").Value = '1.405.81'
$ws.Range("E44").Value = '  +2.72%  '
$ws.Range("D45").Value = '90.56'
$ws.Range("E45").Value = '  -0.44%  '
$ws.Range("E46").Value = '  -4.50%  '
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("D49").Value = '2.91'
$ws.Range("E49").Value = '  +2.27%  '
$ws.Range("D50").Value = '2.205.57'
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("E51").Value = '  -5.75%  '
